# Bill of Materials - DataManagementBoard_V3
# Remove the manual "Backorder" tracking column (H) that was used to flag
# parts on backorder (BACKORDER/X/Backorder/x markers in H2:H17). The
# tracker column is no longer needed, so delete the entire column and
# update the active selection to a cell within the remaining data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column H in its entirety (shifts nothing else - it's the last
# used column) - this also drops the now-unused "BACKORDER"/"X"/
# "Backorder"/"x" shared strings and shrinks the sheet's used range from
# A1:H26 down to A1:F26.
$ws.Columns("H").Delete() | Out-Null

# Move the selection off the deleted column onto a cell still in range.
$ws.Range("C18").Select() | Out-Null
